$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.187.44'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").Value = '3.937.49'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '492.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  -1.10%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000350'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").Value = '4.573.01'
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '3.925.54'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.16%  '
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.78%  '
$ws.Range("D20").Value = '69.353.30'
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '439.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("E23").Value = '  -1.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("E27").Value = '  -4.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '706.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.130'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.463'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +18.48%  '
$ws.Range("D35").Value = '0.0₃0911'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '61.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.152'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0488'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.01'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.143'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.16%  '
$ws.Range("D48").Value = '0.0₆0358'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '
